$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.812.20'
$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D3').Value = '1.584.73'
$ws.Range('E3').Value = '  -2.12%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('E5').Value = '  -1.75%  '
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.481'
$ws.Range('E7').Value = '  -3.68%  '
$ws.Range('E8').Value = '  -1.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0615'
$ws.Range('E9').Value = '  -0.33%  '
$ws.Range('E10').Value = '  -1.75%  '
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('D12').Value = '1.806.06'
$ws.Range('E12').Value = '  -2.01%  '
$ws.Range('D13').Value = '1.583.95'
$ws.Range('E13').Value = '  -1.95%  '
$ws.Range('E14').Value = '  -2.60%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.509'
$ws.Range('E15').Value = '  -2.73%  '
$ws.Range('D16').Value = '25.807.28'
$ws.Range('E16').Value = '  -0.49%  '
$ws.Range('D17').Value = '0.0₃0722'
$ws.Range('E17').Value = '  -2.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '59.78'
$ws.Range('E18').Value = '  -2.91%  '
$ws.Range('E19').Value = '  +0.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '191.19'
$ws.Range('E20').Value = '  -0.11%  '
$ws.Range('E21').Value = '  -1.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.36'
$ws.Range('E22').Value = '  -1.49%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.91'
$ws.Range('E23').Value = '  -2.02%  '
$ws.Range('E24').Value = '  -1.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '141.48'
$ws.Range('E25').Value = '  -1.60%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('E27').Value = '  -1.11%  '
$ws.Range('E28').Value = '  -1.29%  '
$ws.Range('E29').Value = '  -3.08%  '
$ws.Range('E30').Value = '  -5.78%  '
$ws.Range('E31').Value = '  -1.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.10'
$ws.Range('E32').Value = '  -0.60%  '
$ws.Range('E33').Value = '  -2.43%  '
$ws.Range('E34').Value = '  -0.28%  '
$ws.Range('E35').Value = '  -2.19%  '
$ws.Range('D36').Value = '1.098.13'
$ws.Range('E36').Value = '  -2.33%  '
$ws.Range('E37').Value = '  +0.16%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.34'
$ws.Range('E38').Value = '  -1.66%  '
$ws.Range('E39').Value = '  -1.44%  '
$ws.Range('E40').Value = '  -2.20%  '
$ws.Range('E41').Value = '  -7.80%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.816'
$ws.Range('E42').Value = '  +9.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.21'
$ws.Range('E43').Value = '  +2.77%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '93.75'
$ws.Range('E44').Value = '  -4.16%  '
$ws.Range('D45').Value = '1.719.08'
$ws.Range('E45').Value = '  -1.97%  '
$ws.Range('E46').Value = '  +0.21%  '
$ws.Range('E47').Value = '  -0.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '53.12'
$ws.Range('E48').Value = '  -1.78%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0508'
$ws.Range('E50').Value = '  -0.70%  '
$ws.Range('E51').Value = '  -0.07%  '
